# Apply weekly fruit/vegetable price update:
#  - Rows 12 and 13 (the last two data rows) hold the "Primera" and "Segunda"
#    quality records for the latest week. The old values in these two rows
#    are preserved as a new historical snapshot in rows 14 and 15, and rows
#    12/13 are overwritten in place with the new week's data (new date +
#    new volumes/prices).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Snapshot the current (old) rows 12 and 13 values before they are
#    overwritten, so they can be written back out unchanged into the
#    newly appended rows 14 and 15.
$oldRows = @(12, 13)
$savedValues = @{}
foreach ($r in $oldRows) {
    $rowVals = @{}
    for ($c = 1; $c -le 20; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $savedValues[$r] = $rowVals
}

# 2) Write the saved old row 12 / row 13 content into the brand-new rows
#    14 and 15, cell by cell (preserving values and the date column's
#    number format/style).
$newRowMap = @{ 12 = 14; 13 = 15 }
foreach ($oldR in $oldRows) {
    $newR = $newRowMap[$oldR]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($newR, $c).Value = $savedValues[$oldR][$c]
    }
    # Column D (4) carries a special date number format in the source rows;
    # replicate it onto the new row so formatting matches the original cell.
    $ws.Cells.Item($newR, 4).NumberFormat = $ws.Cells.Item($oldR, 4).NumberFormat
}

# 3) Overwrite row 12 ("Primera") with the new week's figures.
$ws.Cells.Item(12, 4).Value = 44609    # D12 Fecha
$ws.Cells.Item(12, 13).Value = 100     # M12 Volumen
$ws.Cells.Item(12, 14).Value = 6500    # N12 Precio minimo
$ws.Cells.Item(12, 15).Value = 7000    # O12 Precio maximo
$ws.Cells.Item(12, 16).Value = 6750    # P12 Precio promedio ponderado
$ws.Cells.Item(12, 19).Value = 3375    # S12 Precio $/Kg

# 4) Overwrite row 13 ("Segunda") with the new week's figures.
$ws.Cells.Item(13, 4).Value = 44609    # D13 Fecha
$ws.Cells.Item(13, 13).Value = 50      # M13 Volumen
$ws.Cells.Item(13, 14).Value = 6000    # N13 Precio minimo
$ws.Cells.Item(13, 15).Value = 6000    # O13 Precio maximo
$ws.Cells.Item(13, 16).Value = 6000    # P13 Precio promedio ponderado
$ws.Cells.Item(13, 19).Value = 3000    # S13 Precio $/Kg

$wb.Save()
